$d = $word.ActiveDocument

# Update the date heading (paragraph 1)
$d.Paragraphs.Item(1).Range.Text = "2025-07-27 Sunday"

# Update each math-problem cell in the table, in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "44+17=",
    "51-18=",
    "72-46=",
    "6+78=",
    "29+33=",
    "16+29=",
    "34+59=",
    "39+6=",
    "79+14=",
    "40-36=",
    "30-12=",
    "74-45=",
    "30-29=",
    "45+27=",
    "95-78=",
    "83-35=",
    "71-3=",
    "28+29=",
    "80-64=",
    "64-17=",
    "51-29=",
    "25+56=",
    "86+6=",
    "91-35=",
    "49+46=",
    "76-37=",
    "68+16=",
    "28+34=",
    "71-38=",
    "7+9=",
    "93-19=",
    "84-37=",
    "39+4=",
    "26+58=",
    "31-15=",
    "27+66=",
    "16+66=",
    "90-36=",
    "73-54=",
    "70-33=",
    "78-59=",
    "29+69=",
    "70-42=",
    "26+45=",
    "94-26=",
    "42-13=",
    "27+39=",
    "82-37=",
    "78-59=",
    "92-77=",
    "18+49=",
    "6+48=",
    "66+6=",
    "19+6=",
    "19+15=",
    "85-9=",
    "40-31=",
    "82-34=",
    "28+5=",
    "9+53=",
    "6+79=",
    "90-82=",
    "65+8=",
    "91-62=",
    "74-35=",
    "48+18=",
    "16+58=",
    "58-49=",
    "83-78=",
    "90-69=",
    "71-52=",
    "59+25=",
    "18+6=",
    "27+69=",
    "9+29=",
    "8+68=",
    "28+55=",
    "85-68=",
    "16+35=",
    "8+46=",
    "67+7=",
    "32-26=",
    "7+74=",
    "26+5=",
    "80-7=",
    "16+29=",
    "45+29=",
    "38+9=",
    "26+59=",
    "30-27=",
    "46+16=",
    "27+39=",
    "74+17=",
    "35+59=",
    "24+48=",
    "31-2=",
    "81-72=",
    "15+26=",
    "24-15=",
    "85+8="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
if (($rows * $cols) -ne $values.Length) {
    Write-Host "WARNING: table size $rows x $cols does not match value count" $values.Length
}

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Done. Updated" $idx "cells."